# Add a new "localdb" command-type column into the hidden '#system' sheet,
# register its data/function list and its named range, and shift every
# defined name / column that sat at or to the right of column N one
# column further to the right (since the new "localdb" column is being
# inserted at column N).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1) Insert a brand-new column at N. Everything that used to live in
#    columns N..AC (macro, mail, number, pdf, rdbms, redis, sms, sound,
#    ssh, step, web, webalert, webcookie, ws, ws.async, xml) shifts one
#    column to the right, landing in O..AD.
# ---------------------------------------------------------------------
$ws.Columns("N").Insert(-4161)

# ---------------------------------------------------------------------
# 2) Populate the new "localdb" column (N) with its header/category name
#    plus the function list for this new command type.
# ---------------------------------------------------------------------
$ws.Range("N1").Value = "localdb"
$ws.Range("N2").Value = "cloneTable(var,source,target)"
$ws.Range("N3").Value = "dropTables(var,tables)"
$ws.Range("N4").Value = "exportCSV(sql,output)"
$ws.Range("N5").Value = "importRecords(var,sourceDb,sql,table)"
$ws.Range("N6").Value = "purge(var)"
$ws.Range("N7").Value = "runSQLs(var,sqls)"

# ---------------------------------------------------------------------
# 3) "localdb" is also added to the master list of command types kept in
#    column A (the "target" list), sorted alphabetically between "json"
#    and "macro". Insert.() in this environment shifts the whole row
#    across every column, so shift column A's values manually, one row
#    at a time, starting from the bottom, to avoid disturbing the other
#    (already-correct) columns.
# ---------------------------------------------------------------------
for ($r = 29; $r -ge 14; $r--) {
    $ws.Cells.Item($r + 1, 1).Value = $ws.Cells.Item($r, 1).Value2
}
$ws.Cells.Item(14, 1).Value = "localdb"

# ---------------------------------------------------------------------
# 4) Update the workbook-level defined names so they keep pointing at
#    the correct (now shifted) ranges, and register the new "localdb"
#    name.
# ---------------------------------------------------------------------
$wb.Names.Item("macro").RefersTo       = "='#system'!`$O`$2:`$O`$4"
$wb.Names.Item("mail").RefersTo        = "='#system'!`$P`$2:`$P`$2"
$wb.Names.Item("number").RefersTo      = "='#system'!`$Q`$2:`$Q`$16"
$wb.Names.Item("pdf").RefersTo         = "='#system'!`$R`$2:`$R`$16"
$wb.Names.Item("rdbms").RefersTo       = "='#system'!`$S`$2:`$S`$7"
$wb.Names.Item("redis").RefersTo       = "='#system'!`$T`$2:`$T`$10"
$wb.Names.Item("sms").RefersTo         = "='#system'!`$U`$2:`$U`$2"
$wb.Names.Item("sound").RefersTo       = "='#system'!`$V`$2:`$V`$5"
$wb.Names.Item("ssh").RefersTo         = "='#system'!`$W`$2:`$W`$9"
$wb.Names.Item("step").RefersTo        = "='#system'!`$X`$2:`$X`$4"
$wb.Names.Item("web").RefersTo         = "='#system'!`$Y`$2:`$Y`$127"
$wb.Names.Item("webalert").RefersTo    = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo   = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo          = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo    = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo         = "='#system'!`$AD`$2:`$AD`$21"
$wb.Names.Item("target").RefersTo      = "='#system'!`$A`$2:`$A`$30"

$wb.Names.Add("localdb", "='#system'!`$N`$2:`$N`$7")
